# Applies the Alvearie -> LinuxForHealth rebranding edit to the
# StructureDefinition-parent-organization-hierarchy-level-code workbook.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/parent-organization-hierarchy-level-code"

# --- Metadata sheet updates ---
# URL
$meta.Range("B2").Value = $newUrl
# Version
$meta.Range("B3").Value = "8.0.0"
# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# The "Extension.url" row (row 5) carries the same canonical URL as a Fixed
# Value (column Q); keep it in sync with the Metadata URL above.
$elements.Range("Q5").Value = $newUrl

# The base "Extension" row (row 2) no longer lists the ele-1/ext-1
# constraint text in its Constraint(s) column (AI); that text now applies
# only to the "Extension.extension" row.
$elements.Range("AI2").Value = ""
